# "feat: remove loading data" -- the OOXML diff actually *adds* two more
# delivery-note rows (6 and 7) to Sheet1, extending the used range from
# A1:G5 to A1:G7. Mirror the existing row layout:
#   A = ID (number)          B = Description (text)
#   C = Amount (number)      D = Price (text, e.g. "40000.00")
#   E = CreatedAt (date/time serial, formatted like existing rows)
#   F = Month (number)       G = TotalAmount (text, e.g. "320000.00")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 -- Product 5
$ws.Cells.Item(6, 1).Value = 16
$ws.Cells.Item(6, 2).Value = "Product 5"
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "40000.00"
$ws.Cells.Item(6, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 5).Value = 45680.98248739923
$ws.Cells.Item(6, 6).Value = 10
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "320000.00"

# Row 7 -- Product 6
$ws.Cells.Item(7, 1).Value = 17
$ws.Cells.Item(7, 2).Value = "Product 6"
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "80000.00"
$ws.Cells.Item(7, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 45681.35013923144
$ws.Cells.Item(7, 6).Value = 10
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "480000.00"
